$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new row at position 6 (shifts existing rows down, "Book & Collections" now at row 7)
$ws.Rows.Item(6).Insert()

# Set new row 6 content - "Books, Newspapers or Magazines"
# Order matches shared string append order seen in target: English, Spanish, Portuguese
$ws.Range("A6").Value = "Books, Newspapers or Magazines"
$ws.Range("C6").Value = "Libros, periódicos o revistas"
$ws.Range("B6").Value = "Livros, Jornais ou Revistas "
$ws.Range("D6").Value = 490199

# Delete the old "Book & Collections" row (now pushed to row 7)
$ws.Rows.Item(7).Delete()

Write-Host "A6:" $ws.Range("A6").Value()
Write-Host "A7:" $ws.Range("A7").Value()
Write-Host "A8:" $ws.Range("A8").Value()
